# milestone1.2 - Alex Hochberg - 3-14-16
#
# Fills in the state-transition rows for states q12 (row13), q19 (row20),
# q22 (row23), q27 (row28), q32 (row33), q37 (row38), q38 (row39) and q46
# (row47): every letter column (B:AA) transitions to q27, except column Q
# (digit) which goes to q1 and column W (whitespace) which goes to q22; the
# '*' column (AR) transitions to q41. Also corrects the q28/qError values
# that had been swapped between the q1 row (row 3) and q2 row (row 2), and
# moves the active selection to AR54.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(13, 20, 23, 28, 33, 38, 39, 47)

foreach ($r in $rows) {
    $ws.Range("B$r`:AA$r").Value = "q27"
    $ws.Range("Q$r").Value = "q1"
    $ws.Range("W$r").Value = "q22"
    $ws.Range("AR$r").Value = "q41"
}

$ws.Range("AS2").Value = "qError"
$ws.Range("AS3").Value = "q28"

[void]$ws.Range("AR54").Select()
